# Update countries & provincias Spain
#
# Refresh the COVID-19 "Pais" sheet with a newer data pull:
#  - header timestamp bumped from 19:05 to 19:35
#  - most touched rows simply get refreshed numbers for the same country
#  - Cabo Verde's case count overtakes Togo's, so they swap ranks/rows
#  - Benin's case count jumps above Islas Feroe/Guadalupe/Gibraltar/
#    Mongolia/Brunei, so it moves up and those five shift down one row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: refresh the "last updated" timestamp -------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 19:35"

# Each entry: row, country, casos totales, nuevos casos, casos activos,
#             recuperados, casos criticos, muertes hoy, muertes
$data = @(
    @(4,   "Estados Unidos", 1675434, 8606, 449591, 1126879, 0, 281, 98964),
    @(11,  "Alemania",       180153,  167,  160300, 11482,   0, 5,   8371),
    @(12,  "Turquia",        156827,  1141, 118694, 33793,   0, 32,  4340),
    @(96,  "Mayotte",        1587,    66,   894,    673,     0, 1,   20),
    @(140, "Cabo Verde",     380,     9,    155,    222,     0, 0,   3),
    @(141, "Togo",           373,     0,    133,    228,     0, 0,   12),
    @(144, "Ruanda",         327,     2,    237,    90,      0, 0,   0),
    @(158, "Benin",          191,     56,   82,     106,     0, 0,   3),
    @(159, "Islas Feroe",    187,     0,    187,    0,       0, 0,   0),
    @(160, "Guadalupe",      156,     0,    115,    28,      0, 0,   13),
    @(161, "Gibraltar",      154,     2,    147,    7,       0, 0,   0),
    @(162, "Mongolia",       141,     0,    32,     109,     0, 0,   0),
    @(163, "Brunei",         141,     0,    136,    4,       0, 0,   1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
